# Automate pulling CBO projections
#
# The workbook tracks three parallel blocks of rows (same id/name pairs, in the
# same order) distinguished by the "source" column (C):
#   rows   2- 55 : source = "current"    (this quarter's freshly pulled data)
#   rows  56-109 : source = "difference" (current - previous)
#   rows 110-163 : source = "previous"   (last quarter's pulled data, kept as a
#                                         baseline for comparison)
#
# Pulling a new round of CBO projections rolls the window forward one quarter:
#   - the data that used to be "current" becomes the new "previous" data, i.e.
#     what was already sitting in the "previous" block (unchanged, already a
#     quarter old) becomes the new "current" block, replacing the old
#     projection numbers there.
#   - "previous" itself is left untouched (it already reflects last quarter's
#     pull).
#   - "difference" is recomputed; right after a fresh pull current == previous
#     so every difference value collapses to 0.
#   - the oldest forecast quarter column ("2024 Q2", column P) is dropped
#     entirely since it is no longer part of the rolling projection window.
#
# Historical sub-rows (id = "historical", columns D:F) never change quarter to
# quarter, only the "projection" sub-rows (columns G:O) are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataCol = 7   # column G
$lastDataCol  = 15  # column O (P is dropped below)

$currentStart  = 2
$currentEnd    = 55
$previousOffset = 108   # previous row = current row + 108
$differenceOffset = 54  # difference row = current row + 54

for ($r = $currentStart; $r -le $currentEnd; $r++) {
    $id = $ws.Cells.Item($r, 1).Value2

    if ($id -eq "projection") {
        $prevRow = $r + $previousOffset
        $diffRow = $r + $differenceOffset

        for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
            # new current = old previous (the window rolls forward)
            $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($prevRow, $c).Value2

            # new difference = new current - new previous = 0, right after the pull
            $ws.Cells.Item($diffRow, $c).Value2 = 0
        }
    }
}

# The oldest forecast quarter ("2024 Q2") is no longer part of the rolling
# window, so drop its column outright.
$ws.Columns("P").Delete()
